$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update class-diagram text cells (UML notation: name : Type) ---
$ws.Range('C5').Value = '+ preberiEhrIdZaVitalneZnake : String '
$ws.Range('K5').Value = 'EHRScape API'
$ws.Range('C6').Value = '+ meritveVitalnihZnakovEHRid : String'
$ws.Range('K6').Value = '# https://rest.ehrscape.com/rest/v1/'
$ws.Range('C7').Value = '+ generiraj()'
$ws.Range('K7').Value = '+ session(username : String, password : String) : String '
$ws.Range('C8').Value = '+ izracunIndexTelesneMase()'
$ws.Range('G8').Value = '<<control>>'
$ws.Range('K8').Value = '+ demographics/EHR/EHRid/parity (sessionID : String, EHRid : String) :  JSON_object '
$ws.Range('C9').Value = '+ graphBMI()'
$ws.Range('G9').Value = 'koda'
$ws.Range('K9').Value = '+ view/EHRid/weight(sessionID : String, EHRid : String) : JSON_array'
$ws.Range('C10').Value = '+ pocistiGraphBMI()'
$ws.Range('G10').Value = '- baseUrl : String '
$ws.Range('K10').Value = '+ view/EHRid/height(sessionID : String, EHRid : String) : JSON_array '
$ws.Range('G11').Value = '- queryUrl : String '
$ws.Range('G12').Value = '- username : String'
$ws.Range('G13').Value = '- password : String '
$ws.Range('G14').Value = '- items : Int_array'
$ws.Range('G15').Value = '- visina : Int_array '
$ws.Range('G16').Value = '- teza : Int_array'
$ws.Range('C17').Value = 'vnos'
$ws.Range('G17').Value = '- dat : Int_array '
$ws.Range('C18').Value = '+ meritevEHRId : String'
$ws.Range('G18').Value = '- gBMI : Int_array '
$ws.Range('C19').Value = '+ meritevVisina : Int'
$ws.Range('G19').Value = '- manj : Int'
$ws.Range('C20').Value = '+ meritevDatum : String '
$ws.Range('G20').Value = '- vec : Int '
$ws.Range('C21').Value = '+ meritevTeza : Int '
$ws.Range('G21').Value = '- scale : Float '
$ws.Range('C22').Value = '+ vnosMeritev()'
$ws.Range('G22').Value = '# getSessionId()'
$ws.Range('G23').Value = '- generirajPodatke(stPacienta : Int)'
$ws.Range('G24').Value = '+ generiraj()'
$ws.Range('G25').Value = '- kreirajNovEhr(ime : String, priimek : String, datumRojstva : String) : String '
$ws.Range('G26').Value = '+ vnosMeritev()'
$ws.Range('G27').Value = '- pridobiStatisticnePodatke()'
$ws.Range('G28').Value = '+ izracunIndexTelesneMase()'
$ws.Range('G29').Value = '+ graphBMI()'
$ws.Range('G30').Value = '- log(text : String)'
$ws.Range('G31').Value = '+ pocistiGraphBMI()'
$ws.Range('G32').Value = '- grafFunkcija()'
$ws.Range('G33').Value = '- Int width(margin : Int)'
$ws.Range('G34').Value = '- Int height(margin : Int)'
$ws.Range('G35').Value = '- sortFunction(a : Int, Int b : Int)'
$ws.Range('G36').Value = '- generateData()'
$ws.Range('G37').Value = '+ $(window).resize(function callback)'
$ws.Range('G38').Value = '+ $(document).ready(function callback)'
$ws.Range('C3').Value = '<<boundary>>'
$ws.Range('C4').Value = 'index'
$ws.Range('K4').Value = '<<boundary>>'
$ws.Range('C16').Value = '<<boundary>>'

# --- Column width adjustments (col G / col 7, col K / col 11) ---
$ws.Columns.Item(7).ColumnWidth = 66.16666666666667
$ws.Columns.Item(11).ColumnWidth = 76.16666666666667

# --- View state: scroll to show row 3 at top, zoom to 100%, select K11 ---
# (ScrollRow/ScrollColumn position the "top-left visible cell" the same way
#  Excel's topLeftCell view attribute does; Zoom is the closest exposed
#  equivalent of the saved zoom percentage.)
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
[void]$ws.Range("K11").Select()
$win.Zoom = 100
